# Swap the data contents of rows 16 and 17 on the active sheet.
# (Row 16 becomes the old row 17's record and vice versa - everything
# except the row number travels with the record: Id, TaxonId, species
# names, coordinates, times, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 gets what used to be in row 17 ---
$ws.Range("A16").Value2 = 111768503
$ws.Range("B16").Value2 = 88966
$ws.Range("D16").Value2 = "NT"
$ws.Range("E16").Value2 = 5754
$ws.Range("F16").Value2 = "Gultoppig fingersvamp"
$ws.Range("G16").Value2 = "Ramaria testaceoflava"
$ws.Range("H16").Value2 = "(Bres.) Corner"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value2 = "20"
$ws.Range("J16").Value2 = "fruktkroppar"
$ws.Range("L16").ClearContents()
$ws.Range("Q16").Value2 = 525545.3455456314
$ws.Range("R16").Value2 = 6727837.787189188
$ws.Range("Z16").Value2 = "15:22"
$ws.Range("AB16").Value2 = "15:22"

# --- Row 17 gets what used to be in row 16 ---
$ws.Range("A17").Value2 = 111768476
$ws.Range("B17").Value2 = 96348
$ws.Range("D17").Value2 = "VU"
$ws.Range("E17").Value2 = 220787
$ws.Range("F17").Value2 = "Knärot"
$ws.Range("G17").Value2 = "Goodyera repens"
$ws.Range("H17").Value2 = "(L.) R. Br."
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value2 = "25"
$ws.Range("J17").Value2 = "plantor/tuvor"
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value2 = ""
$ws.Range("Q17").Value2 = 525546.5036804043
$ws.Range("R17").Value2 = 6727881.884716956
$ws.Range("Z17").Value2 = "15:21"
$ws.Range("AB17").Value2 = "15:21"
